# Add in new ONS data: update the "Job adverts by profession" row with the
# newer ONS Textkernel release (Dec 2022 (13/02/23)) and its updated dataset
# URL, drop the now-unused helper column E, and normalise direct formatting
# on the data cells (only a handful of cells keep bespoke formatting: the
# bold-ish row labels, the ILR source cells, and explicit number formats on
# the "period" columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the "Job adverts by profession" row (row 13) with the latest
#    ONS Textkernel figures / link.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "<a href='https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/labourdemandvolumesbyprofessionandlocalauthorityukjanuary2017todecember2022'>ONS Textkernel</a>"
$ws.Range("C13").Value = "Dec 2022 (13/02/23)"

# ---------------------------------------------------------------------
# 2. Remove the stray, empty helper column E (never held any data).
# ---------------------------------------------------------------------
$ws.Range("E1:E13").Clear() | Out-Null

# ---------------------------------------------------------------------
# 3. Strip the direct formatting that was blanket-applied to (almost)
#    every cell, leaving only the handful of cells that keep bespoke
#    formatting untouched (A2:A6, A13, B5:C5, B6:C6, C13).
# ---------------------------------------------------------------------
$resetRanges = @(
    "A1:D1",
    "B2:B4",
    "D2:D4",
    "D5:D6",
    "A7:D10",
    "A11:B12",
    "B13",
    "D13"
)
foreach ($r in $resetRanges) {
    $ws.Range($r).Style = "Normal"
}

# ---------------------------------------------------------------------
# 4. The "Latest period" column for the Annual Population Survey rows
#    becomes an explicit date-style format (mmm-yy) ...
# ---------------------------------------------------------------------
$dateFmtRange = $ws.Range("C2:C4")
$dateFmtRange.Style = "Normal"
$dateFmtRange.NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------
# ... and the Key Stage 4 / Key Stage 5 period columns become explicit
# text-style formats (@).
# ---------------------------------------------------------------------
$textFmtRanges = @("C11:D11", "C12:D12")
foreach ($r in $textFmtRanges) {
    $rng = $ws.Range($r)
    $rng.Style = "Normal"
    $rng.NumberFormat = "@"
}

# ---------------------------------------------------------------------
# 5. Leave the selection where the user ended up after editing the table.
# ---------------------------------------------------------------------
$ws.Range("B14").Select() | Out-Null
